$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date (serial 45190 = 2023-09-21) that was
# bumped by two days (serial 45192 = 2023-09-23) for every data row (2-342).
for ($r = 2; $r -le 342; $r++) {
    $ws.Cells.Item($r, 3).Value = 45192
}
